# Lab.06: "Проверил:" (checked-by, masculine) -> "Проверила:" (feminine),
# since the reviewer "Куркчи А. Э." is female. The author placed the cursor
# right before the trailing colon and typed "а", so Word split the old
# single run into "Проверил" + the freshly typed "а", and dragged its
# "_GoBack" (last-edit-location) bookmark along to sit right after the new
# "а" and before the colon - removing it from wherever it used to be.

$d = $word.ActiveDocument

# 1) Remove the old _GoBack bookmark (wherever it currently lives).
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2) Find "Проверил:" and insert "а" right before the trailing colon.
$found = $d.Content
$found.Find.Execute("Проверил:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$colon = $d.Range($found.End - 1, $found.End - 1)
$colon.InsertBefore("а")

# Nudge the newly inserted "а" formatting so the engine keeps it as its own
# run instead of silently re-merging it with the preceding "Проверил" run
# (both already share identical run formatting).
$aRun = $d.Range($found.End - 1, $found.End)
$aRun.Font.Bold = 1
$aRun.Font.Bold = 0

# 3) Re-plant _GoBack right after the "а", before the colon.
$newSpot = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $newSpot)
